$d = $word.ActiveDocument

# --- Step 1: paragraph with "Verifying the changes in master" (yellow) ---
# becomes "Verifying the change update" (magenta) with a lastRenderedPageBreak
# marker before the text. Do this one first (it is the later paragraph) so
# that the earlier paragraph's index is unaffected by this replacement.
$pVerify = $d.Paragraphs.Item(22)
$rVerify = $pVerify.Range

$xmlVerify = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
            '<w:p>' +
              '<w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr>' +
              '<w:r>' +
                '<w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:highlight w:val="magenta"/></w:rPr>' +
                '<w:lastRenderedPageBreak/>' +
                '<w:t>Verifying the change update</w:t>' +
              '</w:r>' +
            '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

$rVerify.InsertXML($xmlVerify)

# --- Step 2: the empty paragraph just before it gains a "Changes made for
# testing" (yellow) run, and a brand-new "Let's check now" (green) paragraph
# is inserted right after it, before the "Verifying..." paragraph. ---
$pEmpty = $d.Paragraphs.Item(21)
$rEmpty = $pEmpty.Range

$xmlEmpty = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
            '<w:p>' +
              '<w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr>' +
              '<w:r>' +
                '<w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:highlight w:val="yellow"/></w:rPr>' +
                '<w:t>Changes made for testing</w:t>' +
              '</w:r>' +
            '</w:p>' +
            '<w:p>' +
              '<w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr>' +
              '<w:r>' +
                '<w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:highlight w:val="green"/></w:rPr>' +
                '<w:t>Let’s check now</w:t>' +
              '</w:r>' +
            '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

$rEmpty.InsertXML($xmlEmpty)

Write-Host "edit complete"
